$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Tables{T}" sheet (sheet1.xml): add two new table rows (Invoices /
#    InvoicePositions). Shared-string entries must land in the exact order
#    the original author typed them so the rebuilt sharedStrings.xml keeps
#    the same <si> ordering as the target file.
# ---------------------------------------------------------------------------
$wsTables = $wb.Worksheets.Item("Tables{T}")

# First pass - touches the cells that introduce brand-new shared strings,
# in the order that reproduces indices 213..222.
$wsTables.Range("O18").Value = "Rechnung"
$wsTables.Range("B18").Value = "Invoices"
$wsTables.Range("P18").Value = "Invoice"
$wsTables.Range("B19").Value = "InvoicePositions"
$wsTables.Range("L19").Value = "invoiceId"
$wsTables.Range("O19").Value = "Rechnungposition"
$wsTables.Range("P19").Value = "Invoice position"
$wsTables.Range("S18").Value = "Settings/InvoicePositions"
$wsTables.Range("T18").Value = "invoicePositions"
$wsTables.Range("U18").Value = "Positionen"

# Remaining cells reuse already-known shared strings / plain numbers.
$wsTables.Range("A18").Value = "dbo"
$wsTables.Range("C18").Value = "Configuration"
$wsTables.Range("D18").Value = "Settings"
$wsTables.Range("R18").Value = "x"
$wsTables.Range("V18").Value = "Positions"
$wsTables.Range("AA18").Value = 12
$wsTables.Range("AB18").Value = "Rechnung"
$wsTables.Range("AD18").Value = "x"

$wsTables.Range("A19").Value = "dbo"
$wsTables.Range("C19").Value = "Configuration"
$wsTables.Range("D19").Value = "Settings"
$wsTables.Range("K19").Value = "x"
$wsTables.Range("M19").Value = "Invoices"

# Selection / scroll state for this sheet.
$wsTables.Activate()
$wsTables.Range("A17:F21").Select()

# ---------------------------------------------------------------------------
# 2. "Columns{C}" sheet (sheet2.xml): becomes the active tab, selection
#    moves down the sheet.
# ---------------------------------------------------------------------------
$wsColumns = $wb.Worksheets.Item("Columns{C}")
$wsColumns.Activate()
$wsColumns.Range("A72:I74").Select()

# ---------------------------------------------------------------------------
# 3. "Controllers{T}" sheet (sheet4.xml): was the active tab before, loses
#    that status now that Columns{C} is active (handled automatically by
#    activating Columns{C} last / setting the workbook's active sheet).
# ---------------------------------------------------------------------------

# Make sure the workbook-level active sheet is Columns{C} (matches
# bookViews/workbookView activeTab going from 3 to 1).
$wsColumns.Activate()
